$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Spaces"
$ws.Range("F2").Value = "     "
$ws.Range("F3").Value = "         "
$ws.Range("F4").Value = " "

$ws.Range("G10").Select() | Out-Null
